# Add two more "test" rows (row 3 and row 4) to every results sheet, and
# refresh the existing row 2 values, so all five algorithm sheets report the
# same three test cases (commit: "ajuste para fazer os mesmos testes para os
# algoritmos.").

$wb = $excel.ActiveWorkbook

# Per-sheet data for rows 2, 3 and 4.
# Columns: B=Origem, C=Destino, D=Caminho, E=Quantidade de nós expandidos,
#          F=Fator de ramificação médio, G=Tempo, H=Memória Alocada
# NOTE: the PS parser here doesn't understand scientific-notation numeric
# literals (e.g. "7.05e-05"), so every "Tempo" value below is written out in
# plain decimal - same IEEE-754 double, just without the exponent.
$data = @{
    "BFS" = @(
        @{ B=2; C=4; D="[2, 5, 4]"; E=6; F=1.166666666666667;  G=0.0001971721649169922;      H=0 },
        @{ B=3; C=3; D="[3]";       E=1; F=0;                   G=0.0000705718994140625;      H=0 },
        @{ B=9; C=9; D="[9]";       E=1; F=0;                   G=0.00007319450378417969;     H=0 }
    )
    "DFS" = @(
        @{ B=2; C=4; D="[2, 3, 6, 5, 4]"; E=6; F=0; G=0.00007939338684082031; H=0 },
        @{ B=3; C=3; D="[3]";             E=1; F=0; G=0.0001015663146972656;  H=0 },
        @{ B=9; C=9; D="[9]";             E=1; F=0; G=0.00009512901306152344; H=0 }
    )
    "BCU" = @(
        @{ B=2; C=4; D="[2, 5, 4]"; E=8; F=0.7; G=0.00009799003601074219; H=0 },
        @{ B=3; C=3; D="[3]";       E=1; F=0;   G=0.00006723403930664062; H=0 },
        @{ B=9; C=9; D="[9]";       E=1; F=0;   G=0.00006365776062011719; H=0 }
    )
    "A_Estrela_Euclidiano" = @(
        @{ B=2; C=4; D="[2, 5, 4]"; E=3; F=2.692307692307693; G=0.0001680850982666016;  H=0 },
        @{ B=3; C=3; D="[3]";       E=1; F=2.692307692307693; G=0.00007367134094238281; H=0 },
        @{ B=9; C=9; D="[9]";       E=1; F=2.692307692307693; G=0.00007104873657226562; H=0 }
    )
    "A_Estrela_Haversiano" = @(
        @{ B=2; C=4; D="[2, 5, 4]"; E=13; F=2.692307692307693; G=0.006186723709106445;   H=0 },
        @{ B=3; C=3; D="[3]";       E=1;  F=2.692307692307693; G=0.00008153915405273438; H=0 },
        @{ B=9; C=9; D="[9]";       E=1;  F=2.692307692307693; G=0.0001342296600341797;  H=0 }
    )
}

foreach ($sheetName in $data.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $data[$sheetName]

    for ($i = 0; $i -lt $rows.Count; $i++) {
        $r = $i + 2
        $row = $rows[$i]

        # Column A ("Índice") is a text "1" on every data row, same as the
        # existing row 2 - copy it down instead of re-typing so it keeps the
        # same (string) cell type/style rather than becoming a number.
        if ($r -gt 2) {
            $ws.Cells.Item(2, 1).Copy($ws.Cells.Item($r, 1))
        }

        $ws.Cells.Item($r, 2).Value = $row.B
        $ws.Cells.Item($r, 3).Value = $row.C
        $ws.Cells.Item($r, 4).Value = $row.D
        $ws.Cells.Item($r, 5).Value = $row.E
        $ws.Cells.Item($r, 6).Value = $row.F
        $ws.Cells.Item($r, 7).Value = $row.G
        $ws.Cells.Item($r, 8).Value = $row.H
    }
}
